# Updated cryptos list (price / volume(1h) refresh) - GitHub Actions style data sync.
# D-column "Price" cells are forced to Text ("@") before the write so values that
# parse as numbers (e.g. "585.19", "1.00", "3.10") keep their exact textual
# formatting instead of being coerced into numeric cells; the style is then reset
# back to "Normal" so no stray number-format style sticks to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '68.028.13'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -0.21%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.610.19'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -1.46%  '
$ws.Range("E4").Value = '  -0.10%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '585.19'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.51%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '193.04'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.26%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '3.602.15'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +0.19%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.677'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -2.85%  '
$ws.Range("E11").Value = '  -0.77%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '55.25'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -3.00%  '
$ws.Range("E13").Value = '  +6.37%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '9.98'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -2.63%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '4.182.37'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -1.62%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '3.608.44'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -1.57%  '
$ws.Range("E17").Value = '  -0.32%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '12.50'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.68%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '67.921.73'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("E20").Value = '  -2.11%  '
$ws.Range("E21").Value = '  -2.44%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '404.13'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.11%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '13.46'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +23.46%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '4.24'
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '85.93'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -2.59%  '
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '12.56'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +0.46%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '3.93'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +6.63%  '
$ws.Range("E29").Value = '  +0.65%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '8.09'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +12.96%  '
$ws.Range("E31").Value = '  -2.08%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '31.52'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -1.20%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '683.08'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +12.20%  '
$ws.Range("E34").Value = '  -0.22%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.118'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +1.78%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '64.45'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -5.35%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '42.65'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -3.76%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.422'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +8.08%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("E40").Value = '  +2.06%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '2.96'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +17.80%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '3.196.84'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +15.44%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '3.11'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +7.77%  '
$ws.Range("E44").Value = '  -0.84%  '
$ws.Range("E45").Value = '  -0.11%  '
$ws.Range("E46").Value = '  -0.69%  '
$ws.Range("E47").Value = '  -2.80%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '3.10'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -3.46%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '8.77'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -1.29%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '142.44'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -0.66%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '2.56'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.72%  '
